$wb = $excel.ActiveWorkbook

# 1. About sheet: correct the text referencing the model country (India -> Mexico)
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A12").Value = "For the Mexico model, the desired output units are:"

# 2. CDCF-PMpPDOU sheet: correct conversion factor (replace formula with corrected literal value)
$wsPM = $wb.Worksheets.Item("CDCF-PMpPDOU")
$wsPM.Range("B2").Value = 621372736649.80676

# 3. CDCF-FTMpFDOU sheet: correct conversion factor (replace formula with corrected literal value)
$wsFTM = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$wsFTM.Range("B2").Value = 621372736649.80676
